# Rename the first sheet from "Sheet1" to "commands"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "commands"

# Add the new row of data (use case + command) at the bottom of the commands sheet
$ws1.Range("A30").Value = "How to verify logs of a multi-container pod?"
$ws1.Range("B30").Value = "kubectl logs -f webapp --all-containers=true"

# Copy the style from the row above (A29) to A30 so it matches existing data rows
$null = $ws1.Range("A29").Copy()
$null = $ws1.Range("A30").PasteSpecial(-4122) # xlPasteFormats

# Select A2 on the commands sheet and make it the active sheet/tab
$null = $ws1.Activate()
$null = $ws1.Range("A2").Select()

# Make the imperative sheet not the tab-selected one (commands should be active)
$ws2 = $wb.Worksheets.Item("imperative")
$null = $ws2.Range("B3").Select()
$null = $ws1.Activate()
